$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-09-26 Friday" "2025-09-27 Saturday"

Replace-Text "237×7=" "481×7="
Replace-Text "621×5=" "982×9="
Replace-Text "659×6=" "241×4="
Replace-Text "523×4=" "959×3="
Replace-Text "206×9=" "428×9="
Replace-Text "795×4=" "368×8="
Replace-Text "773×5=" "548×7="
Replace-Text "969×2=" "738×6="
Replace-Text "198×5=" "644×9="
Replace-Text "646×5=" "505×5="
Replace-Text "219×3=" "530×2="
Replace-Text "106×7=" "585×8="
Replace-Text "548×8=" "313×9="
Replace-Text "826×4=" "198×2="
Replace-Text "526×7=" "989×6="
Replace-Text "832×8=" "240×3="
Replace-Text "844×3=" "386×2="
Replace-Text "607×6=" "775×4="
Replace-Text "639×4=" "838×7="
Replace-Text "880×6=" "549×3="
Replace-Text "865×3=" "678×6="
Replace-Text "597×7=" "431×6="
Replace-Text "183×6=" "903×8="
Replace-Text "494×9=" "569×7="
Replace-Text "455×6=" "588×2="
